$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I0 and IF in row 1, matching the style of the existing header cell H1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Values for I2:I33 and J2:J33 as described in the diff
$iValues = @(4,6,8,8,12,9,2,3,9,8,7,3,4,5,8,7,8,6,9,3,6,7,6,5,8,5,8,6,7,5,4,4)
$jValues = @(4,7,8,8,12,9,2,3,9,8,8,3,5,5,8,7,9,6,9,5,6,7,6,6,8,5,8,6,8,5,4,4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
